$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8258154988288879
$ws.Range("B1").Value = 1.694992423057556
$ws.Range("C1").Value = 6.030293464660645
$ws.Range("D1").Value = 1.952608346939087
$ws.Range("E1").Value = 1.171127319335938
